$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.418.44"
$ws.Range("E2").Value = "  +2.84%  "
$ws.Range("D3").Value = "2.015.62"
$ws.Range("E3").Value = "  +6.19%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'246.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.52%  "
$ws.Range("D6").Value = "'0.660"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").Value = "'45.43"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.62%  "
$ws.Range("D9").Value = "'0.363"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.48%  "
$ws.Range("D10").Value = "'56.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("D11").Value = "'0.0717"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.92%  "
$ws.Range("D12").Value = "'0.0988"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("D13").Value = "'14.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.58%  "
$ws.Range("D14").Value = "2.300.79"
$ws.Range("E14").Value = "  +5.92%  "
$ws.Range("D15").Value = "'0.802"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.58%  "
$ws.Range("D16").Value = "2.019.72"
$ws.Range("E16").Value = "  +6.50%  "
$ws.Range("D17").Value = "'4.89"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.00%  "
$ws.Range("D18").Value = "36.490.47"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").Value = "'70.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.66%  "
$ws.Range("D20").Value = "0.0₃0814"
$ws.Range("E20").Value = "  -1.65%  "
$ws.Range("D21").Value = "'12.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.41%  "
$ws.Range("D22").Value = "'233.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.50%  "
$ws.Range("D23").Value = "'4.98"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -5.39%  "
$ws.Range("D24").Value = "'1.01"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.58%  "
$ws.Range("D25").Value = "'2.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.22%  "
$ws.Range("D26").Value = "'161.70"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.14%  "
$ws.Range("D27").Value = "'19.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.63%  "
$ws.Range("D28").Value = "'1.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -8.71%  "
$ws.Range("D29").Value = "'8.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.41%  "
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("D31").Value = "'21.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +57.75%  "
$ws.Range("D32").Value = "'4.37"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.39%  "
$ws.Range("D33").Value = "'0.0583"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.81%  "
$ws.Range("D34").Value = "'0.999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").Value = "'0.0859"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +20.98%  "
$ws.Range("D37").Value = "'4.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.91%  "
$ws.Range("D38").Value = "'2.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +11.61%  "
$ws.Range("D39").Value = "'0.844"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("D40").Value = "'1.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.01%  "
$ws.Range("D41").Value = "'97.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.79%  "
$ws.Range("D42").Value = "'0.0215"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.38%  "
$ws.Range("D43").Value = "'16.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.60%  "
$ws.Range("D44").Value = "'1.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.79%  "
$ws.Range("D45").Value = "'2.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.24%  "
$ws.Range("D46").Value = "1.308.51"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").Value = "'0.0812"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").Value = "'2.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("D49").Value = "'2.24"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.31%  "
$ws.Range("D50").Value = "2.204.26"
$ws.Range("E50").Value = "  +6.30%  "
$ws.Range("D51").Value = "'3.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +13.42%  "
